$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 16.558025
$ws.Range("N2").Value = 33.11605
$ws.Range("O2").Value = 0.2047872600336892
$ws.Range("P2").Value = 0.155585835247525
$ws.Range("Q2").Value = 0.4840655818625
$ws.Range("R2").Value = 1.93626232745
$ws.Range("S2").Value = 0.09069546744010251
$ws.Range("T2").Value = 0.05389283104016977
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("M3").Value = 20.31779433333333
$ws.Range("O3").Value = 0.2512875437409564
$ws.Range("P3").Value = 0.2863712008291233
$ws.Range("Q3").Value = 0.5939805584378333
$ws.Range("R3").Value = 3.563883350627
$ws.Range("S3").Value = 0.1112893508986448
$ws.Range("T3").Value = 0.09919511449420315
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("M4").Value = 10.405091
$ws.Range("N4").Value = 31.215273
$ws.Range("O4").Value = 0.1286886616182304
$ws.Range("P4").Value = 0.1466556042216543
$ws.Range("Q4").Value = 0.3041876328395
$ws.Range("R4").Value = 1.825125797037
$ws.Range("S4").Value = 0.05699318560044471
$ws.Range("T4").Value = 0.05079951967887997
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("M5").Value = 13.1587975
$ws.Range("N5").Value = 26.317595
$ws.Range("O5").Value = 0.1627461056112162
$ws.Range("P5").Value = 0.1236453320906656
$ws.Range("Q5").Value = 0.38469086551375
$ws.Range("R5").Value = 1.538763462055
$ws.Range("S5").Value = 0.07207642760607935
$ws.Range("T5").Value = 0.0428290723295386
$ws.Range("G6").Value = 0.0292345
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.4428765120700495
$ws.Range("J6").Value = 0.346386487911515
$ws.Range("M6").Value = 7.303315666666667
$ws.Range("N6").Value = 21.909947
$ws.Range("O6").Value = 0.09032635260170116
$ws.Range("P6").Value = 0.1029373190408882
$ws.Range("Q6").Value = 0.2135087818571667
$ws.Range("R6").Value = 1.281052691143
$ws.Range("S6").Value = 0.04000341998825085
$ws.Range("T6").Value = 0.03565609641760037
$ws.Range("G7").Value = 0.0292345
$ws.Range("H7").Value = 0.058469
$ws.Range("I7").Value = 0.4428765120700495
$ws.Range("J7").Value = 0.346386487911515
$ws.Range("M7").Value = 13.11173766666667
$ws.Range("N7").Value = 39.335213
$ws.Range("O7").Value = 0.1621640763942067
$ws.Range("P7").Value = 0.1848047085701436
$ws.Range("Q7").Value = 0.3833150948161667
$ws.Range("R7").Value = 2.299890568897
$ws.Range("S7").Value = 0.0718186605365273
$ws.Range("T7").Value = 0.06401385395112309
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.653613512088485
$ws.Range("M8").Value = 16.558025
$ws.Range("N8").Value = 33.11605
$ws.Range("O8").Value = 0.2047872600336892
$ws.Range("P8").Value = 0.155585835247525
$ws.Range("Q8").Value = 0.6089379274
$ws.Range("R8").Value = 3.6536275644
$ws.Range("S8").Value = 0.1140917925935867
$ws.Range("T8").Value = 0.1016930042073552
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.653613512088485
$ws.Range("M9").Value = 20.31779433333333
$ws.Range("O9").Value = 0.2512875437409564
$ws.Range("P9").Value = 0.2863712008291233
$ws.Range("Q9").Value = 0.7472072044026665
$ws.Range("R9").Value = 6.724864839623999
$ws.Range("S9").Value = 0.1399981928423116
$ws.Range("T9").Value = 0.1871760863349201
$ws.Range("I10").Value = 0.5571234879299505
$ws.Range("J10").Value = 0.653613512088485
$ws.Range("M10").Value = 10.405091
$ws.Range("N10").Value = 31.215273
$ws.Range("O10").Value = 0.1286886616182304
$ws.Range("P10").Value = 0.1466556042216543
$ws.Range("Q10").Value = 0.3826576266159999
$ws.Range("R10").Value = 3.443918639544
$ws.Range("S10").Value = 0.07169547601778564
$ws.Range("T10").Value = 0.09585608454277429
$ws.Range("I11").Value = 0.5571234879299505
$ws.Range("J11").Value = 0.653613512088485
$ws.Range("M11").Value = 13.1587975
$ws.Range("N11").Value = 26.317595
$ws.Range("O11").Value = 0.1627461056112162
$ws.Range("P11").Value = 0.1236453320906656
$ws.Range("Q11").Value = 0.4839279368599999
$ws.Range("R11").Value = 2.90356762116
$ws.Range("S11").Value = 0.09066967800513688
$ws.Range("T11").Value = 0.080816259761127
$ws.Range("I12").Value = 0.5571234879299505
$ws.Range("J12").Value = 0.653613512088485
$ws.Range("M12").Value = 7.303315666666667
$ws.Range("N12").Value = 21.909947
$ws.Range("O12").Value = 0.09032635260170116
$ws.Range("P12").Value = 0.1029373190408882
$ws.Range("Q12").Value = 0.2685867369573333
$ws.Range("R12").Value = 2.417280632616
$ws.Range("S12").Value = 0.05032293261345031
$ws.Range("T12").Value = 0.06728122262328778
$ws.Range("I13").Value = 0.5571234879299505
$ws.Range("J13").Value = 0.653613512088485
$ws.Range("M13").Value = 13.11173766666667
$ws.Range("N13").Value = 39.335213
$ws.Range("O13").Value = 0.1621640763942067
$ws.Range("P13").Value = 0.1848047085701436
$ws.Range("Q13").Value = 0.4821972644293334
$ws.Range("R13").Value = 4.339775379864
$ws.Range("S13").Value = 0.09034541585767937
$ws.Range("T13").Value = 0.1207908546190205
